$d = $word.ActiveDocument

# --- Locate the two paragraphs that need to be merged -----------------
# Paragraph A ends with "... cointegration based model"
# Paragraph B is "for forecasting high frequency financial time series". "
$paraA = $null
$paraB = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*cointegration based model*") {
        $paraA = $p
    }
    if ($paraA -ne $null -and $p.Range.Text -like "for forecasting high frequency financial time series*" -and $paraB -eq $null) {
        $paraB = $p
    }
}

$paraAEnd = $paraA.Range.End          # position right after the paragraph mark of A
$paraAEndBeforeMark = $paraAEnd - 1   # position right before the paragraph mark of A (end of its text)

# --- Grab a same-formatted space character to splice in as its own run
# (re-uses the space between "Please" and "find" at the very start of
# paragraph A, which already carries sz=28/szCs=28/lang=en-GB formatting)
$spaceSrc = $d.Range($paraA.Range.Start + 6, $paraA.Range.Start + 7)
$spaceFT = $spaceSrc.FormattedText

# --- Splice in the space first, right before paragraph A's own mark ---
$insPoint1 = $d.Range($paraAEndBeforeMark, $paraAEndBeforeMark)
$insPoint1.FormattedText = $spaceFT

# --- Now grab paragraph B's text (minus its own trailing paragraph
# mark) *fresh*, after the mutation above, and splice it in right after
# the space (again right before paragraph A's, now shifted, mark) -----
$paraBFT = $d.Range($paraB.Range.Start, $paraB.Range.End - 1).FormattedText
$paraAEndBeforeMark2 = $paraA.Range.End - 1
$insPoint2 = $d.Range($paraAEndBeforeMark2, $paraAEndBeforeMark2)
$insPoint2.FormattedText = $paraBFT

# --- Remove the now-duplicated original paragraph B (content + mark) --
$paraB.Range.Delete()

# --- Move the "_GoBack" bookmark to the end of the merged paragraph ---
# (Word only keeps a single "_GoBack" bookmark; adding one removes any
# pre-existing one automatically, matching the diff which relocates it
# from next to "x-" to the end of this paragraph.)
$newParaAEnd = $paraA.Range.End
$bookmarkRange = $d.Range($newParaAEnd - 1, $newParaAEnd - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

Write-Output "done"
